# Initial protected feature imputation
# Adds a new "Protected feature" task row to the "Data exploration" sheet
# (owner: Burke, status: 1st pass, with a note about bio-data-only accuracy),
# and updates the active sheet/selection to reflect the edited location.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Data exploration")
$ws2 = $wb.Worksheets.Item("Modeling")

# New row 12 on the "Data exploration" sheet
$ws1.Range("B12").Value = "   Protected feature"
$ws1.Range("C12").Value = "Burke"
$ws1.Range("E12").Value = "1st pass"
$ws1.Range("F12").Value = "With bio data only, about 70% accuracy is achieved across different classifiers."

# Restore Modeling sheet's prior selection, then make Data exploration the
# active sheet again with its new selection on the freshly-added row.
$ws2.Activate()
$ws2.Range("I5").Select()

$ws1.Activate()
$ws1.Range("D10").Select()
